$wb = $excel.ActiveWorkbook

# The JUnit sheet (first sheet, tab name "JUnit") gets two new test rows
# appended beneath the existing data (rows 1-9 -> now rows 1-11).
$ws = $wb.Worksheets.Item("JUnit")

# Row 10: testAddUserRoleSuccess
$ws.Cells.Item(10, 2).Value = 5
$ws.Cells.Item(10, 3).Value = "testAddUserRoleSuccess"
$ws.Cells.Item(10, 5).Value = "addUserRole"
$ws.Cells.Item(10, 6).Value = "Success added new User Role to database."

# Row 11: testGetAllUserRoleSuccess
$ws.Cells.Item(11, 2).Value = 6
$ws.Cells.Item(11, 3).Value = "testGetAllUserRoleSuccess"
$ws.Cells.Item(11, 5).Value = "getAllRecords"
$ws.Cells.Item(11, 6).Value = "Success get all user roles from the database."

# Move the tracked selection to just past the new last row, matching the
# author's final cursor position after adding the rows.
$ws.Range("F12").Select()
